$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price values
$ws.Range("I2").Value = 260
$ws.Range("I6").Value = 820
$ws.Range("I7").Value = 400

# Update selected cell / active selection
$ws.Range("I3").Select()
